# "Key- Value Data Store.pptx" update
#
#   1) Shape "Rectangle 3": close the unterminated bracket -
#        "[Either with default or Custom File Path"
#      becomes
#        "[Either with default or Custom File Path]"
#
#   2) Shape "Rectangle 4": fix the "Tack" -> "Task" typo. In the canonical
#      edit this was typed as a correction over the selected word/prefix,
#      which leaves the text split across two runs:
#        "Keep Task " + "of keys that"
#      (previously a single run "Keep Tack of keys that").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# --- 1) Rectangle 3 --------------------------------------------------------
$rect3 = Get-ShapeByName $s "Rectangle 3"
if ($rect3 -ne $null) {
    $tr3 = $rect3.TextFrame.TextRange
    $oldFragment = "[Either with default or Custom File Path"
    $newFragment = "[Either with default or Custom File Path]"
    $pos0 = $tr3.Text.IndexOf($oldFragment)
    if ($pos0 -ge 0) {
        $sub3 = $tr3.Characters($pos0 + 1, $oldFragment.Length)
        $sub3.Text = $newFragment
    }
}

# --- 2) Rectangle 4 --------------------------------------------------------
$rect4 = Get-ShapeByName $s "Rectangle 4"
if ($rect4 -ne $null) {
    $tr4 = $rect4.TextFrame.TextRange
    $oldPrefix = "Keep Tack "
    $newPrefix = "Keep Task "
    $pos0 = $tr4.Text.IndexOf($oldPrefix)
    if ($pos0 -ge 0) {
        $sub4 = $tr4.Characters($pos0 + 1, $oldPrefix.Length)
        $sub4.Text = $newPrefix
    }
}
